$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values update
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: move value from D2 to C2 with updated figure, clear D2
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 42.117013549239076

# Row 3: clear B3 and C3 (values removed)
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update the selection to match the new extent used in the workbook
$ws.Range("B1:E3").Select()
